# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B (TB), C (d2S), D (K), E (IP) values per row (F = Win is unchanged).
# G (sum) = B + C + D + E, recomputed to match.
$data = @(
    @{ Row = 2; B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 },
    @{ Row = 3; B = 1.505614041169197;  C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 6.741336633845642 },
    @{ Row = 4; B = 0.7287194209349384; C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 5.964442013611383 },
    @{ Row = 5; B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 },
    @{ Row = 6; B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 },
    @{ Row = 7; B = 0.7287194209349384; C = 1.65323645889881;  D = 3.082599426703578;  E = 6.48142807727062;   G = 11.94598338380795 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 7).Value = $entry.G
}
